$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A2: Student ID 212254 -> 211877, kept as text (reuse E2's plain-text style) ---
$ws.Range("H1").Formula = '=TEXT(211877,"0")'
$ws.Range("H1").Copy()
$ws.Range("A2").PasteSpecial(-4163)   # xlPasteValues: write the literal text, no number re-parsing
$ws.Range("H1").Clear()
$ws.Range("E2").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats: restore the original (unchanged) style

# --- C2: Log Date 09/09/2025 -> 14/10/2025 (stays plain text; style untouched) ---
$ws.Range("C2").Value = "14/10/2025"

# --- D2: Log Time becomes a literal text value "10:30:00" instead of the numeric time serial ---
$ws.Range("H1").Formula = '="10:30:00"'
$ws.Range("H1").Copy()
$ws.Range("D2").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("H1").Clear()
$ws.Range("E2").Copy()
$ws.Range("D2").PasteSpecial(-4122)   # xlPasteFormats: adopt the plain-text style (like the other text cells)
